# Commit: "Fixed issues with 81RF protective element. Changed default xls
# parameters to disable 81x protections. Added goose messages for DER's cb's"
#
# On sheet "relays", for every relay row (2-11) the 81RF rate-of-change
# protection settings are reset to their "disabled" defaults in both the
# first (T:V) and duplicate/second (AH:AJ) 81RF column groups:
#   81RFRP  [Hz/Sec] :   2 -> 100
#   81RFDFP [Hz]      :  57 ->  10
#   81RF Trip Delay   : 0.01 -> 0.1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

for ($row = 2; $row -le 11; $row++) {
    $ws.Range("T$row").Value  = 100
    $ws.Range("U$row").Value  = 10
    $ws.Range("V$row").Value  = 0.1

    $ws.Range("AH$row").Value = 100
    $ws.Range("AI$row").Value = 10
    $ws.Range("AJ$row").Value = 0.1
}

# Leave the sheet focused with the last-edited (second) 81RF block selected,
# matching the reviewer's final on-screen selection.
$ws.Activate()
$ws.Range("AH2:AJ11").Select()
